$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "REPORT"
$ws.Range("AK3").Select()
